# Updated cryptos list on Wed Dec 27 02:08:24 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'42.591.48"
$ws.Range('D2').ClearFormats()
$ws.Range('E2').Value = '  -2.43%  '
$ws.Range('D3').Value = "'2.230.65"
$ws.Range('D3').ClearFormats()
$ws.Range('E3').Value = '  -1.95%  '
$ws.Range('D4').Value = "'1.01"
$ws.Range('D4').ClearFormats()
$ws.Range('E4').Value = '  +0.29%  '
$ws.Range('D5').Value = "'111.52"
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  -6.74%  '
$ws.Range('D6').Value = "'291.64"
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  +9.39%  '
$ws.Range('D7').Value = "'0.624"
$ws.Range('D7').ClearFormats()
$ws.Range('E7').Value = '  -2.88%  '
$ws.Range('E8').Value = '  -0.24%  '
$ws.Range('E9').Value = '  -2.36%  '
$ws.Range('D10').Value = "'44.09"
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  -7.49%  '
$ws.Range('D11').Value = "'0.0918"
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  -2.68%  '
$ws.Range('D12').Value = "'54.77"
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  +1.73%  '
$ws.Range('D13').Value = "'8.75"
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  -7.55%  '
$ws.Range('D14').Value = "'1.04"
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  +15.64%  '
$ws.Range('D15').Value = "'0.103"
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  -2.96%  '
$ws.Range('D16').Value = "'15.00"
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '  -4.33%  '
$ws.Range('D17').Value = "'2.563.46"
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '  -2.07%  '
$ws.Range('D18').Value = "'2.231.23"
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  -1.96%  '
$ws.Range('D19').Value = "'42.497.93"
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  -2.56%  '
$ws.Range('D20').Value = "'7.22"
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  +4.33%  '
$ws.Range('E21').Value = '  -3.72%  '
$ws.Range('D22').Value = "'73.03"
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  +1.06%  '
$ws.Range('D23').Value = "'3.40"
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  +17.01%  '
$ws.Range('D24').Value = "'2.41"
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  +0.00%  '
$ws.Range('D25').Value = "'230.34"
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  -1.81%  '
$ws.Range('D26').Value = "'9.15"
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  -3.64%  '
$ws.Range('D27').Value = "'1.00"
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  -1.74%  '
$ws.Range('D28').Value = "'11.54"
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  -4.12%  '
$ws.Range('E29').Value = '  -1.44%  '
$ws.Range('D30').Value = "'37.77"
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  -10.28%  '
$ws.Range('D31').Value = "'173.95"
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  -0.14%  '
$ws.Range('D32').Value = "'3.18"
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  -4.70%  '
$ws.Range('D33').Value = "'21.01"
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  -2.22%  '
$ws.Range('E34').Value = '  -3.32%  '
$ws.Range('D35').Value = "'5.68"
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  -0.86%  '
$ws.Range('D36').Value = "'5.09"
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  +10.38%  '
$ws.Range('D37').Value = "'4.30"
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  +2.05%  '
$ws.Range('E38').Value = '  -2.22%  '
$ws.Range('D39').Value = "'0.0373"
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  -3.10%  '
$ws.Range('E40').Value = '  -3.62%  '
$ws.Range('D41').Value = "'75.05"
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  +3.10%  '
$ws.Range('D42').Value = "'2.43"
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  -4.41%  '
$ws.Range('D43').Value = "'0.233"
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  -2.46%  '
$ws.Range('E44').Value = '  -0.05%  '
$ws.Range('D45').Value = "'12.55"
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  -9.10%  '
$ws.Range('E46').Value = '  -6.03%  '
$ws.Range('E47').Value = '  -5.84%  '
$ws.Range('B48').Value = 'Stacks'
$ws.Range('C48').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D48').Value = "'1.76"
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  +13.49%  '
$ws.Range('B49').Value = 'TrustWalletToken'
$ws.Range('C49').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D49').Value = "'1.30"
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  +2.75%  '
$ws.Range('D50').Value = "'103.03"
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  -0.87%  '
$ws.Range('D51').Value = "'8.46"
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  -0.96%  '
